$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Police spoke with grade 10 student who stated she believes a black van followed her as she walked to school from her residence on three occasions. All three incidents occurred in the morning. She provided limited information due to language barrier however stated that the driver was a "male with a beard" and that the van was newer, black in color and had a partial plate that includes ABC. No physical or verbal interaction or attempts by the driver; he only followed her for a few blocks. Parents were contacted and a safety plan has been put in place including for victim to not walk to school alone in the mornings for the time being.'
